# Auto-generated PowerShell COM-interop script
$wb = $excel.ActiveWorkbook

# 1) Insert the new worksheet 'A06 vie saint gregoire' right after sheet 5
$afterSheet = $wb.Worksheets.Item(5)
$greg = $wb.Worksheets.Add($null, $afterSheet)
$greg.Name = 'A06 vie saint gregoire'

# Header row
$greg.Cells.Item(1,1).Value = 'line_n'
$greg.Cells.Item(1,2).Value = 'prev_line'
$greg.Cells.Item(1,3).Value = 'line'
$greg.Cells.Item(1,4).Value = 'next_line'
$headerRange = $greg.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Data rows
$greg.Cells.Item(2,1).Value = 85
$greg.Cells.Item(2,2).Value = 'Puis a dit a son frere, sans point de demourer:'
$greg.Cells.Item(2,3).Value = '“Portez ce panyer en ung batel en la mer,'
$greg.Cells.Item(2,4).Value = 'Et puis si le laissiez tout a par luy aller.'
$greg.Cells.Item(3,1).Value = 104
$greg.Cells.Item(3,2).Value = 'A l''abbé le porterent, si l''ont mis a raison:'
$greg.Cells.Item(3,3).Value = '“Sire, nous ne povons point de poisson peschier,'
$greg.Cells.Item(3,4).Value = 'Mais nous avons trouvé sur la mer ung panyer;'
$greg.Cells.Item(4,1).Value = 114
$greg.Cells.Item(4,2).Value = 'L''abbé regarde es tables, si dist ceste raison:'
$greg.Cells.Item(4,3).Value = '“Il est de noble lieu; il portera mon nom.'
$greg.Cells.Item(4,4).Value = 'Se Dieu plaist il sera encore moult preudom.”'
$greg.Cells.Item(5,1).Value = 227
$greg.Cells.Item(5,2).Value = 'La chose lui conta, lequel fut en esmoy:'
$greg.Cells.Item(5,3).Value = '“Et non pourtant a lui je me combateray,'
$greg.Cells.Item(5,4).Value = 'Et se je le puis vaincre, foy que doy saint Martin,'
$greg.Cells.Item(6,1).Value = 258
$greg.Cells.Item(6,2).Value = 'Puis est venu au conte, si lui dist a hault ton:'
$greg.Cells.Item(6,3).Value = '“Mon cheual avez mort, c''est trop grant traÿson.'
$greg.Cells.Item(6,4).Value = 'Mais bien tost en arez le vostre guerredon.”'
$greg.Cells.Item(7,1).Value = 336
$greg.Cells.Item(7,2).Value = 'Quant elle vit les tables, elle dit haultement:'
$greg.Cells.Item(7,3).Value = '“Mon filz m''a espousee! De dueil le cueur me fent.'
$greg.Cells.Item(7,4).Value = 'Avecques moy sept ans a esté tellement.”'
$greg.Cells.Item(8,1).Value = 361
$greg.Cells.Item(8,2).Value = 'Adonc recommenca la povre gent a braire.'
$greg.Cells.Item(8,3).Value = '“Demourez avec nous, gentilz homs debonnaire.'
$greg.Cells.Item(8,4).Value = 'Quant en voulez aller, forment nous doit desplaire.'
$greg.Cells.Item(9,1).Value = 368
$greg.Cells.Item(9,2).Value = 'Son aumosnier appelle, si lui dist doulcement:'
$greg.Cells.Item(9,3).Value = '“Donnez leur a chacun ung gros tournois d''argent.'
$greg.Cells.Item(9,4).Value = 'Je m''en yray toudis, trop me font de tourment;'
$greg.Cells.Item(10,1).Value = 404
$greg.Cells.Item(10,2).Value = 'Lors l''oste respondi, qui n''estoit point vilain:'
$greg.Cells.Item(10,3).Value = '“Pres de cy a vne ysle que vous verrez demain.'
$greg.Cells.Item(10,4).Value = 'Il ya vne roche, je vous dy de certain,'
$greg.Cells.Item(11,1).Value = 424
$greg.Cells.Item(11,2).Value = 'Puis a dit a son hoste moult debonairement:'
$greg.Cells.Item(11,3).Value = '“Or refermez la roche, tost et delivrement,'
$greg.Cells.Item(11,4).Value = 'Et men bailliez la clef, car je la vous demand.”'
$greg.Cells.Item(12,1).Value = 472
$greg.Cells.Item(12,2).Value = 'Quant l''oste la regarde, si dist ceste raison:'
$greg.Cells.Item(12,3).Value = '“C''est la clef de la roche, sans variatïon,'
$greg.Cells.Item(12,4).Value = 'Ou j''enfermay Gregoire, qui tant estoit preudom.”'
$greg.Cells.Item(13,1).Value = 478
$greg.Cells.Item(13,2).Value = 'Et il leur print a dire, sans gaires demourer:'
$greg.Cells.Item(13,3).Value = '“Sept ans a que il fut en la roche enfermé.'
$greg.Cells.Item(13,4).Value = 'Je croy qu''il est piecha du siecle trespassé.”'
$greg.Cells.Item(14,1).Value = 492
$greg.Cells.Item(14,2).Value = 'Et ilz lui respondirent a vne voix trestous:'
$greg.Cells.Item(14,3).Value = '“Vous en venrez, beau sire, maintenant avec nous,'
$greg.Cells.Item(14,4).Value = 'En la cité de Romme; povoir arez sur tous.'
$greg.Cells.Item(15,1).Value = 499
$greg.Cells.Item(15,2).Value = 'Son hoste lui a dit, sans longue demourance:'
$greg.Cells.Item(15,3).Value = '“Beau sire veez la, n''en soyez en doubtance.”'
$greg.Cells.Item(15,4).Value = 'Quant il la vit si dist: “Dieu, qu''avez grant puissace!'
$greg.Cells.Item(16,1).Value = 548
$greg.Cells.Item(16,2).Value = 'Devant lui est venue, haultement lui escrie:'
$greg.Cells.Item(16,3).Value = '“Pere, plus grant pecheresse ne fut onc mais ouÿe.”'
$greg.Cells.Item(16,4).Value = 'Saint Gregoire l''apelle, et lui dist: “Doulce amye,'

# 2) Shift the numeric prefix of all subsequent sheet names by +1
#    (these sheets are now one position further back, at index+1,
#    because of the newly inserted sheet above; we find them by
#    their current (old) name so the order of operations doesn't matter)
($wb.Worksheets.Item('A06 saint jean evangeliste')).Name = 'A07 saint jean evangeliste'
($wb.Worksheets.Item('A07 vie saint jean paulus')).Name = 'A08 vie saint jean paulus'
($wb.Worksheets.Item('A08 vie glorieux confesseur')).Name = 'A09 vie glorieux confesseur'
($wb.Worksheets.Item('A09 vie saint leu')).Name = 'A10 vie saint leu'
($wb.Worksheets.Item('A10 poines enfer')).Name = 'A11 poines enfer'
($wb.Worksheets.Item('A11 vie saint sebastien')).Name = 'A12 vie saint sebastien'
($wb.Worksheets.Item('A12 miracle saint servais')).Name = 'A13 miracle saint servais'
($wb.Worksheets.Item('A13 vie seint thibault')).Name = 'A14 vie seint thibault'
($wb.Worksheets.Item('A16 guillaume angleterre')).Name = 'A17 guillaume angleterre'
($wb.Worksheets.Item('A17 robert deable')).Name = 'A18 robert deable'
($wb.Worksheets.Item('A18 richart sans peour')).Name = 'A19 richart sans peour'
($wb.Worksheets.Item('A19 elegy troyes')).Name = 'A20 elegy troyes'
($wb.Worksheets.Item('A20 vieillards tués')).Name = 'A21 vieillards tués'
($wb.Worksheets.Item('A21 mauvais riche homme')).Name = 'A22 mauvais riche homme'
($wb.Worksheets.Item('A22 jeu des dez')).Name = 'A23 jeu des dez'
($wb.Worksheets.Item('A23 roy avoit amie')).Name = 'A24 roy avoit amie'
($wb.Worksheets.Item('A25 quatre sereurs')).Name = 'A26 quatre sereurs'

# 3) Fix the quotation marks on the 'guillaume angleterre' sheet (now 'A17 guillaume angleterre'),
#    row with line_n=182, column D (next_line)
$gaSheet = $wb.Worksheets.Item('A17 guillaume angleterre')
$gaSheet.Range("D4").Value = '“Hé! roy,” dist la roïne, “vos amours fausses sont;'

